# Auto-generated script applying scheduled market-data refresh to Siren_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for the affected leve rows
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching the upstream commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 12821584
$ws.Range("I6").Value = 33333690
$ws.Range("J6").Value = 1517.9375
$ws.Range("K6").Value = 100001070
$ws.Range("L6").Value = 4553.8125
$ws.Range("M6").Value = -100000958
$ws.Range("N6").Value = -4777.8125
$ws.Range("H76").Value = 4660.769
$ws.Range("I76").Value = 3719.4
$ws.Range("K76").Value = 3719.4
$ws.Range("M76").Value = -3404.4
$ws.Range("H79").Value = 4660.769
$ws.Range("I79").Value = 3719.4
$ws.Range("K79").Value = 3719.4
$ws.Range("M79").Value = -2627.4
$ws.Range("H110").Value = 468033340
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H135").Value = 5315.2104
$ws.Range("I135").Value = 5602.5
$ws.Range("K135").Value = 50422.5
$ws.Range("M135").Value = -47887.5
$ws.Range("H138").Value = 226419.06
$ws.Range("I138").Value = 659118.9
$ws.Range("K138").Value = 1977356.7
$ws.Range("M138").Value = -1972216.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8318.299999999999
$ws.Range("I45").Value = 8030.5
$ws.Range("K45").Value = 8030.5
$ws.Range("M45").Value = -7653.5
$ws.Range("H102").Value = 14120.25
$ws.Range("I102").Value = 27596.25
$ws.Range("K102").Value = 27596.25
$ws.Range("M102").Value = -25974.25
$ws.Range("H132").Value = 2472.743
$ws.Range("I132").Value = 814.5
$ws.Range("K132").Value = 2443.5
$ws.Range("M132").Value = 86.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2246.2104
$ws.Range("I20").Value = 1575.6428
$ws.Range("K20").Value = 1575.6428
$ws.Range("M20").Value = -1328.6428
$ws.Range("H94").Value = 3263.6365
$ws.Range("I94").Value = 2064.7273
$ws.Range("J94").Value = 5661.4546
$ws.Range("K94").Value = 2064.7273
$ws.Range("L94").Value = 5661.4546
$ws.Range("M94").Value = -1613.7273
$ws.Range("N94").Value = -6563.4546
$ws.Range("H105").Value = 5760.3667
$ws.Range("I105").Value = 6474.1
$ws.Range("K105").Value = 6474.1
$ws.Range("M105").Value = -4727.1
$ws.Range("H107").Value = 6380.6665
$ws.Range("I107").Value = 6803.25
$ws.Range("K107").Value = 6803.25
$ws.Range("M107").Value = -4883.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 180.375
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 191.85715
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 191.85715
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -417.85715
$ws.Range("H86").Value = 4686.2188
$ws.Range("I86").Value = 4320.7
$ws.Range("K86").Value = 4320.7
$ws.Range("M86").Value = -3197.7
$ws.Range("H89").Value = 4686.2188
$ws.Range("I89").Value = 4320.7
$ws.Range("K89").Value = 21603.5
$ws.Range("M89").Value = -15987.5
$ws.Range("H130").Value = 65000
$ws.Range("J130").Value = 65000
$ws.Range("L130").Value = 65000
$ws.Range("N130").Value = -75040
$ws.Range("H134").Value = 4601.1113
$ws.Range("J134").Value = 6713.857
$ws.Range("L134").Value = 20141.571
$ws.Range("N134").Value = -25211.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 36.5
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 42.75
$ws.Range("K2").Value = 204
$ws.Range("L2").Value = 256.5
$ws.Range("M2").Value = -91
$ws.Range("N2").Value = -482.5
$ws.Range("H34").Value = 1668.1818
$ws.Range("I34").Value = 1621.5385
$ws.Range("J34").Value = 1735.5555
$ws.Range("K34").Value = 4864.6155
$ws.Range("L34").Value = 5206.666499999999
$ws.Range("M34").Value = -4780.6155
$ws.Range("N34").Value = -5374.666499999999
$ws.Range("H75").Value = 241.66667
$ws.Range("I75").Value = 100
$ws.Range("J75").Value = 270
$ws.Range("K75").Value = 300
$ws.Range("L75").Value = 810
$ws.Range("M75").Value = 698
$ws.Range("N75").Value = -2806
$ws.Range("H78").Value = 241.66667
$ws.Range("I78").Value = 100
$ws.Range("J78").Value = 270
$ws.Range("K78").Value = 900
$ws.Range("L78").Value = 2430
$ws.Range("M78").Value = 4092
$ws.Range("N78").Value = -12414
$ws.Range("H141").Value = 2632.9
$ws.Range("I141").Value = 2369.889
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 7109.667
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -1929.667
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 215673.25
$ws.Range("I20").Value = 838669
$ws.Range("K20").Value = 838669
$ws.Range("M20").Value = -838424
$ws.Range("H24").Value = 895117.3
$ws.Range("H70").Value = 10289.857
$ws.Range("I70").Value = 8180.6665
$ws.Range("K70").Value = 8180.6665
$ws.Range("M70").Value = -7910.6665
$ws.Range("H73").Value = 10289.857
$ws.Range("I73").Value = 8180.6665
$ws.Range("K73").Value = 8180.6665
$ws.Range("M73").Value = -7244.6665
$ws.Range("H132").Value = 2391.9167
$ws.Range("I132").Value = 2524.1428
$ws.Range("K132").Value = 7572.428400000001
$ws.Range("M132").Value = -5042.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 21466.812
$ws.Range("I7").Value = 26205.834
$ws.Range("K7").Value = 26205.834
$ws.Range("M7").Value = -26093.834
$ws.Range("H22").Value = 1312.5
$ws.Range("I22").Value = 625
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 625
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -330
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1312.5
$ws.Range("I27").Value = 625
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 625
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -518
$ws.Range("N27").Value = -2214
$ws.Range("H55").Value = 1328.2941
$ws.Range("I55").Value = 558.1
$ws.Range("K55").Value = 558.1
$ws.Range("M55").Value = -385.1
$ws.Range("H126").Value = 21466.812
$ws.Range("I126").Value = 26205.834
$ws.Range("K126").Value = 78617.50199999999
$ws.Range("M126").Value = -76147.50199999999
$ws.Range("H132").Value = 393946.62
$ws.Range("I132").Value = 533437.75
$ws.Range("K132").Value = 1600313.25
$ws.Range("M132").Value = -1597783.25
$ws.Range("H136").Value = 4780.1953
$ws.Range("I136").Value = 3262
$ws.Range("K136").Value = 9786
$ws.Range("M136").Value = -7236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H96").Value = 7144525
$ws.Range("I96").Value = 12500997
$ws.Range("K96").Value = 12500997
$ws.Range("M96").Value = -12499624
$ws.Range("H107").Value = 25121.46
$ws.Range("I107").Value = 3222.6667
$ws.Range("J107").Value = 43891.855
$ws.Range("K107").Value = 9668.000100000001
$ws.Range("L107").Value = 131675.565
$ws.Range("M107").Value = -7748.000100000001
$ws.Range("N107").Value = -135515.565
$ws.Range("H132").Value = 6794.0933
$ws.Range("I132").Value = 7613.4907
$ws.Range("K132").Value = 22840.4721
$ws.Range("M132").Value = -20310.4721
$ws.Range("H136").Value = 336298.44
$ws.Range("I136").Value = 396269.9
$ws.Range("J136").Value = 2171.7144
$ws.Range("K136").Value = 1188809.7
$ws.Range("L136").Value = 6515.1432
$ws.Range("M136").Value = -1186259.7
$ws.Range("N136").Value = -11615.1432
